$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$pairs = @(
    ,("49+4=","37+55=")
    ,("55+18=","56+18=")
    ,("29+9=","72-28=")
    ,("36+16=","39+34=")
    ,("72-46=","65+7=")
    ,("7+39=","83-16=")
    ,("40-17=","15+76=")
    ,("91-79=","65-6=")
    ,("93-19=","82-63=")
    ,("15+46=","80-22=")
    ,("69+6=","61-9=")
    ,("34-26=","18+64=")
    ,("84+8=","80-26=")
    ,("28+6=","86+7=")
    ,("36+57=","71-29=")
    ,("28+43=","90-15=")
    ,("83-26=","47+48=")
    ,("74-8=","71-33=")
    ,("73+19=","49+49=")
    ,("46+49=","76+6=")
    ,("84-69=","82-37=")
    ,("73-14=","27+58=")
    ,("60-36=","71-56=")
    ,("16+25=","65-47=")
    ,("61-15=","85-29=")
    ,("9+59=","60-1=")
    ,("81-6=","46+47=")
    ,("8+87=","72-53=")
    ,("95-89=","22-17=")
    ,("90-75=","6+16=")
    ,("19+68=","46+8=")
    ,("8+28=","7+14=")
    ,("45+48=","14+77=")
    ,("62-3=","30-18=")
    ,("45+47=","69+6=")
    ,("36+47=","80-21=")
    ,("59+26=","7+85=")
    ,("61-28=","90-51=")
    ,("41-25=","69+16=")
    ,("81-18=","71-64=")
    ,("77-48=","27+35=")
    ,("41-3=","32+59=")
    ,("17+14=","71-19=")
    ,("43-18=","47-8=")
    ,("53-37=","13+58=")
    ,("29+67=","84-48=")
    ,("35+29=","22+9=")
    ,("28+28=","70-5=")
    ,("9+65=","81-33=")
    ,("54-49=","80-14=")
    ,("43+38=","65-28=")
    ,("76-48=","83-14=")
    ,("28+4=","91-89=")
    ,("90-35=","90-77=")
    ,("82-17=","19+47=")
    ,("37-19=","94-25=")
    ,("88-29=","27+44=")
    ,("6+29=","30-25=")
    ,("61-58=","6+47=")
    ,("60-8=","29+3=")
    ,("89+5=","8+77=")
    ,("69+29=","18+49=")
    ,("9+69=","90-79=")
    ,("94-19=","97-48=")
    ,("53-27=","86-69=")
    ,("49+26=","46-7=")
    ,("39+35=","85+7=")
    ,("46+7=","90-12=")
    ,("40-28=","92-75=")
    ,("51-44=","44-9=")
    ,("91-12=","68+14=")
    ,("56-8=","76+15=")
    ,("64+18=","91-66=")
    ,("93-8=","21-14=")
    ,("86+5=","12-9=")
    ,("72-8=","90-8=")
    ,("17+47=","47+9=")
    ,("18+23=","18+63=")
    ,("40-26=","24+58=")
    ,("33+9=","69+4=")
    ,("73-26=","95-6=")
    ,("25+38=","97-88=")
    ,("91-6=","47+44=")
    ,("56-28=","91-38=")
    ,("93-8=","84-28=")
    ,("33-19=","6+86=")
    ,("79+17=","51-9=")
    ,("41-3=","82-33=")
    ,("68+7=","93-45=")
    ,("74-25=","81-63=")
    ,("51-5=","13+39=")
    ,("15+27=","64-19=")
    ,("15+38=","27+18=")
    ,("45-28=","88-39=")
    ,("25+49=","6+58=")
    ,("28+8=","28+48=")
    ,("94-59=","83-9=")
    ,("38+49=","18+5=")
    ,("12-6=","58+39=")
    ,("47+5=","60-2=")
)
$cols = 5
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = [int]([math]::Floor($i / $cols)) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $pairs[$i][1]
}
Write-Output "done"
